# Generate Report for Handoff
#
# - Status moves from "In Translation" to "Ready for handoff" (Overview!E2,
#   Overview!F2, zh-cn!C2, de-de!C2).
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" /
#   "Latest Handback DateTime" timestamps are refreshed to the new handoff
#   run's generation time.
# - The Status / date columns get a bit wider so the new "Ready for
#   handoff" text isn't truncated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ---------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps ----------------------------------------
# Overview!G2 and de-de!H2 shared the same "Latest Handoff" timestamp
# before the edit (the de-de handoff landed in the same run as the
# Overview roll-up); zh-cn's handoff ran a few seconds earlier/later and
# keeps its own timestamp.
$overview.Range("G2").Value = "2016-08-29 19:04:30"
$dede.Range("H2").Value     = "2016-08-29 19:04:30"
$zhcn.Range("H2").Value     = "2016-08-29 19:04:25"

# --- Widen the Status / date columns so the longer text fits -------------
$overview.Range("E1").ColumnWidth = 16.38265482584637
$overview.Range("F1").ColumnWidth = 16.38265482584637
$zhcn.Range("C1").ColumnWidth     = 16.38265482584637
$dede.Range("C1").ColumnWidth     = 16.38265482584637
